$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.660.40'
$ws.Range("E2").Value = '  +0.00%  '

$ws.Range("D3").Value = '3.324.82'
$ws.Range("E3").Value = '  +0.82%  '

$ws.Range("E4").Value = '  +0.17%  '

$ws.Range("D5").Value = '''581.02'
$ws.Range("E5").Value = '  +0.87%  '

$ws.Range("D6").Value = '''175.37'
$ws.Range("E6").Value = '  -1.31%  '

$ws.Range("E7").Value = '  -0.15%  '

$ws.Range("E8").Value = '  +0.93%  '

$ws.Range("D9").Value = '3.321.20'
$ws.Range("E9").Value = '  +0.86%  '

$ws.Range("D10").Value = '''0.180'
$ws.Range("E10").Value = '  +3.07%  '

$ws.Range("D11").Value = '''0.579'
$ws.Range("E11").Value = '  +1.11%  '

$ws.Range("D12").Value = '''46.43'
$ws.Range("E12").Value = '  +1.70%  '

$ws.Range("E13").Value = '  +0.57%  '

$ws.Range("D14").Value = '''704.53'
$ws.Range("E14").Value = '  +0.27%  '

$ws.Range("D15").Value = '3.870.82'
$ws.Range("E15").Value = '  +1.02%  '

$ws.Range("D16").Value = '''8.43'
$ws.Range("E16").Value = '  +0.91%  '

$ws.Range("D17").Value = '67.694.88'
$ws.Range("E17").Value = '  -0.13%  '

$ws.Range("E18").Value = '  -0.85%  '

$ws.Range("D19").Value = '3.334.90'
$ws.Range("E19").Value = '  +0.90%  '

$ws.Range("D20").Value = '''17.35'
$ws.Range("E20").Value = '  -0.44%  '

$ws.Range("E21").Value = '  +2.05%  '

$ws.Range("E22").Value = '  +0.07%  '

$ws.Range("D23").Value = '''5.36'
$ws.Range("E23").Value = '  +3.82%  '

$ws.Range("D24").Value = '''16.92'
$ws.Range("E24").Value = '  +0.64%  '

$ws.Range("D25").Value = '''98.57'
$ws.Range("E25").Value = '  -0.04%  '

$ws.Range("D26").Value = '''3.87'
$ws.Range("E26").Value = '  -1.33%  '

$ws.Range("E27").Value = '  -1.71%  '

$ws.Range("D28").Value = '''9.44'
$ws.Range("E28").Value = '  +1.19%  '

$ws.Range("D29").Value = '''33.06'
$ws.Range("E29").Value = '  +0.08%  '

$ws.Range("E30").Value = '  +0.86%  '

$ws.Range("D31").Value = '''7.09'
$ws.Range("E31").Value = '  +5.97%  '

$ws.Range("D32").Value = '''568.86'
$ws.Range("E32").Value = '  -1.57%  '

$ws.Range("D33").Value = '''10.97'
$ws.Range("E33").Value = '  +1.27%  '

$ws.Range("E34").Value = '  +1.55%  '

$ws.Range("D35").Value = '''57.50'
$ws.Range("E35").Value = '  +3.99%  '

$ws.Range("D37").Value = '3.698.75'
$ws.Range("E37").Value = '  -5.36%  '

$ws.Range("D38").Value = '''3.31'
$ws.Range("E38").Value = '  -1.26%  '

$ws.Range("D39").Value = '''34.09'
$ws.Range("E39").Value = '  +6.35%  '

$ws.Range("D40").Value = '''0.130'
$ws.Range("E40").Value = '  +0.74%  '

$ws.Range("E41").Value = '  +1.41%  '

$ws.Range("E42").Value = '  +0.81%  '

$ws.Range("E43").Value = '  -0.79%  '

$ws.Range("E44").Value = '  +1.56%  '

$ws.Range("D45").Value = '''3.31'
$ws.Range("E45").Value = '  -0.60%  '

$ws.Range("E46").Value = '  -1.02%  '

$ws.Range("D47").Value = '''2.67'
$ws.Range("E47").Value = '  +5.32%  '

$ws.Range("E48").Value = '  +0.22%  '

$ws.Range("E49").Value = '  -0.42%  '

$ws.Range("E50").Value = '  -5.80%  '

$ws.Range("D51").Value = '''128.86'
$ws.Range("E51").Value = '  +0.45%  '
